# Auto-generated Excel COM-interop script to update the cryptos list
# matching the GitHub Actions scheduled refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 43 and 44 swap coin identity (Filecoin <-> RenderToken) together
# with their link, price and volume figures.
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

# Price (column D) updates. Cells whose new text parses as a plain
# number need an explicit Text number format first, otherwise Excel
# would silently store them as a numeric value instead of text.
$ws.Range("D2").Value = "59.218.90"
$ws.Range("D3").Value = "2.494.07"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.77"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.26"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D9").Value = "2.518.12"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0995"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.26"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.335"
$ws.Range("D14").Value = "2.942.40"
$ws.Range("D15").Value = "58.884.50"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.47"
$ws.Range("D18").Value = "2.503.96"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.70"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.04"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.15"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.67"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.412"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.52"
$ws.Range("D29").Value = "0.0₃0769"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.66"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.33"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.23"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.53"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.67"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.788"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "280.83"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.21"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.49"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "131.85"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.595"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0932"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0510"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.16"
$ws.Range("D51").Value = "1.761.24"

# Volume(1h) (column E) updates.
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("E6").Value = "  +5.11%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("E9").Value = "  +3.32%  "
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("E14").Value = "  +2.84%  "
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("E17").Value = "  +3.85%  "
$ws.Range("E18").Value = "  +2.81%  "
$ws.Range("E19").Value = "  +3.03%  "
$ws.Range("E20").Value = "  +3.99%  "
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  +9.39%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  +3.71%  "
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("E28").Value = "  +4.43%  "
$ws.Range("E29").Value = "  +7.51%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("E32").Value = "  +4.63%  "
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("E39").Value = "  +5.48%  "
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("E43").Value = "  +7.38%  "
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("E45").Value = "  +10.16%  "
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("E48").Value = "  +5.92%  "
$ws.Range("E49").Value = "  +4.70%  "
$ws.Range("E50").Value = "  +4.46%  "
$ws.Range("E51").Value = "  +3.48%  "
